$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76
$ws.Cells.Item($row, 1).Value = "2025-04-29 11:28:30"
$ws.Cells.Item($row, 2).Value = 225
